# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same 8-column fund-holdings layout we
#    need) right after itself, rename the copy to "2022-Q1" and overwrite
#    its cell contents with the new quarter's fund-holding rows (it already
#    lands in the correct tab position, directly before "总计").
# 2. Insert a new row at the top of "总计"'s data (row 2) for the 2022-Q1
#    summary line and shift the existing 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as TEXT
# (so numeric-looking strings like "012027" or "6.76" don't get silently
# coerced into numbers and lose their formatting / leading zeros).
# Uses a scratch cell formatted as Text, then copies only the *value*
# (which keeps the text type) onto the real target - this avoids leaving
# any new style on the target cell itself.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Item(1).Range("ZZ1")
$scratch.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# =======================================================================
# Step 1: create the "2022-Q1" sheet
# =======================================================================
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)
$newSheet = $wb.Worksheets.Item($srcSheet.Index + 1)
$newSheet.Name = "2022-Q1"

# The source sheet only had 3 data rows (r2:r4); we need a 4th (r5).
# Copy the formatting of row 4 down onto row 5 first.
$newSheet.Range("A4:H4").Copy()
$newSheet.Range("A5:H5").PasteSpecial(-4122)  # xlPasteFormats

# ---- Row 1: headers ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Row 2 ----
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "160642"
$newSheet.Range("C2").Value = "鹏华增瑞灵活配置混合(LOF)"
Set-TextValue $newSheet.Range("D2") "6.76"
Set-TextValue $newSheet.Range("E2") "91.34"
Set-TextValue $newSheet.Range("F2") "3.45"
Set-TextValue $newSheet.Range("G2") "0.2332"
$newSheet.Range("H2").Value = 10

# ---- Row 3 ----
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "012027"
$newSheet.Range("C3").Value = "光大保德信安阳一年持有期混合型证券投资基金A"
Set-TextValue $newSheet.Range("D3") "15.22"
Set-TextValue $newSheet.Range("E3") "22.05"
Set-TextValue $newSheet.Range("F3") "0.96"
Set-TextValue $newSheet.Range("G3") "0.1461"
$newSheet.Range("H3").Value = 3

# ---- Row 4 ----
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "012028"
$newSheet.Range("C4").Value = "光大保德信安阳一年持有期混合型证券投资基金C"
Set-TextValue $newSheet.Range("D4") "7.68"
Set-TextValue $newSheet.Range("E4") "22.05"
Set-TextValue $newSheet.Range("F4") "0.96"
Set-TextValue $newSheet.Range("G4") "0.0737"
$newSheet.Range("H4").Value = 3

# ---- Row 5 ----
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "005444"
$newSheet.Range("C5").Value = "光大保德信多策略精选18个月定期开放灵活配置混合"
Set-TextValue $newSheet.Range("D5") "1.00"
Set-TextValue $newSheet.Range("E5") "29.09"
Set-TextValue $newSheet.Range("F5") "2.02"
Set-TextValue $newSheet.Range("G5") "0.0202"
$newSheet.Range("H5").Value = 6

# =======================================================================
# Step 2: update the "总计" sheet - insert the 2022-Q1 summary row
# =======================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row inherited the header row's formatting; strip it back
# to the plain (no style) look the other data rows use.
$totalSheet.Range("A2:D2").ClearFormats()

# Column A on the data rows uses the same style as the header (s=2) -
# copy that formatting from the header row's A cell equivalent (A3, which
# kept its original style after the insert).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
Set-TextValue $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.47

# The pre-existing rows kept their old row-index values (0, 1) after the
# shift; renumber them to match their new position (1, 2).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Clean up the scratch cell used for forcing text values.
$scratch.Clear()

# Copying a sheet makes the new copy the active tab; restore the
# original active sheet ("2021-Q3", the first tab) so it stays untouched.
$wb.Worksheets.Item(1).Activate()
